# Rerunning slender tuna models
# Adds a new "SLT" (slender tuna) row of model-output results to both
# worksheets: "optimal models" and "relative importance".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "optimal models"
$ws2 = $wb.Worksheets.Item(2)   # "relative importance"

# --- Sheet 1: "optimal models" -------------------------------------------
$ws1.Range("A13").Value = "SLT"
$ws1.Range("B13").Value = 2
$ws1.Range("C13").Value = 0.5
$ws1.Range("D13").Value = 0.009
$ws1.Range("E13").Value = 1300
$ws1.Range("F13").Value = 0.9971
$ws1.Range("G13").Value = 0.98146
$ws1.Range("H13").Value = 0.9854
$ws1.Range("I13").Value = 0.030286

# --- Sheet 2: "relative importance" ---------------------------------------
$ws2.Range("A13").Value = "SLT"
$ws2.Range("B13").Value = 5.2984744
$ws2.Range("C13").Value = 17.114367
$ws2.Range("D13").Value = 7.4964237
$ws2.Range("E13").Value = 1.3195522
$ws2.Range("F13").Value = 7.6811202
$ws2.Range("G13").Value = 6.6096603
$ws2.Range("H13").Value = 0.8076748
$ws2.Range("I13").Value = 0.6772072
$ws2.Range("J13").Value = 11.439719
$ws2.Range("K13").Value = 2.8052522
$ws2.Range("L13").Value = 6.9927041
$ws2.Range("M13").Value = 3.3238051
$ws2.Range("N13").Value = 4.9016924
$ws2.Range("O13").Value = 18.8945761
$ws2.Range("P13").Value = 1.3678342
$ws2.Range("Q13").Value = 2.2101546
$ws2.Range("R13").Value = 1.0597825

# --- Update the saved selection on each sheet to the next empty row ------
$ws1.Range("A14").Select() | Out-Null
$ws2.Range("A14").Select() | Out-Null

# Restore sheet1 ("optimal models") as the active/visible tab, matching
# the original workbook (it was the selected tab before the edit too).
$ws1.Activate()
